$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24-78 down to 25-79.
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with its data.
$ws.Range("A24").Value = 4
$ws.Range("B24").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C24").Value = "Los Lagos"
$ws.Range("D24").Value = 44498
$ws.Range("E24").Value = 10
$ws.Range("F24").Value = 100112022
$ws.Range("G24").Value = "Arveja Verde"
$ws.Range("H24").Value = "Perfection"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 90
$ws.Range("K24").Value = 22000
$ws.Range("L24").Value = 22000
$ws.Range("M24").Value = 22000
$ws.Range("N24").Value = "$/malla 25 kilos"
$ws.Range("O24").Value = "Provincia de Huasco"
$ws.Range("P24").Value = 880
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"

# Keep the date-formatted number format used by the rest of column D.
$ws.Range("D24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
